$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    'Wild moment police officer is knocked out after man allegedly resisting arrest sparks late-night brawl, as three are charged'
    'Funny reactions as Yesu Wa Tongaren allegedly seeks bursary for kids'
    'Crisis in Imo APC as stalwarts resign from partisan politics'
    'Marion County Sheriff''s Office rescues 13 animals and makes arrest in animal cruelty case'
    'Warning to anyone who drinks fizzy drinks - doctor claims they can trigger 2 health issues'
    'Over 30,000 arrested in nationwide anti-drug and crime operations'
    'Two arrested for embezzling funds using Little Hearts project name'
    'PSNI make arrest after female pedestrian ''struck by vehicle'''
    'Read why our founder is building andwemet: A story of rejection and redefining love…'
    'Taxi driver arrested after crash with ''Welcome to Bolton'' road sign'
    'Woman ‘pushed to floor & sexually assaulted’ by entrance to historic castle as man, 20, arrested & cops search for BMW'
    'Secret probes into ‘inappropriate sexual behaviour’ in Archbishop of Wales’ home diocese'
    'News24 | Self-confessed gay imam who was allegedly going to marry lesbian couple gunned down'
    'Out on bail, man held for vehicle theft'
    'Xabi Alonso’s ‘Real Madrid clause’ update as Carlo Ancelotti faces fresh challenge'
    '''''Rogue landlord'' who had 18 people living in three-bedroom home including babies is fined £37k'
    '‘The new wave of crime’: Car part theft is plaguing Melbourne, here’s why'
    'Casa Orsola controversy: Barcelona city government buys the building, creating moral hazard'
    'Meenakshi Natarajan takes charge of Telangana Cong affairs ahead of local polls'
    'Greece cracks down on suspected agricultural subsidy fraud'
    'Burton´s Udoka Godwin-Malife suffers alleged racist...'
    'Here Come the HELOCs: Mortgages, Housing-Debt-to-Income-Ratio, Serious Delinquencies & Foreclosures in Q4 2024'
    'Property management company FirstPort suspended by trade body after breach of rules'
    'Michael Magruder: Financial Advisor''s Costly Loan Scandal - Investment & Securities Fraud Lawyer'
    'Central Vermont man with arrest warrant trespasses, found week later police say'
    'Mark Hughes´ first home game in charge of bottom side...'
    'Slidell man arrested for intentionally parking truck on train track, causing crash'
    'ROSEN, NATIONAL TRIAL LAWYERS, Encourages Walgreens Boots Alliance, Inc. Investors to Secure Counsel Before Important Deadline in Securities Class Action - WBA'
    'Cabcharge review ordered as hospital blasts ‘utterly disgraceful, unethical and unacceptable’ conduct'
    'B.C. drug rehab facility faces closure over abuse, safety violation allegations'
    'TOP RANKED ROSEN LAW FIRM Encourages Light & Wonder, Inc. Investors to Inquire About Securities Class Action Investigation - LNW'
    'Family tributes as rail firm fined £4m for deaths'
    'Kansas woman charged with wire fraud over $260k handbag purchases'
    'Man arrested after showing gun to neighbor during argument in Little Torch Key: MCSO'
    'Evans in charge of four-way fight to win Rally Sweden'
    '''''Mass Theft'': Thousands of Artists Call for AI Art Auction to be Cancelled'
    'Man arrested after intentionally parking truck on train tracks during Slidell parade'
    'They Were VIP Gamblers With Betting Problems. Now, They’re Suing DraftKings'
    'Neuroscience doctors to decide on ending strike after director’s resignation notice'
    'News24 | Principal of top Pretoria school to face six charges of misconduct at disciplinary hearing'
    'Winnipeg Police Arrest Suspect After Random Assault Downtown'
    '1000s of Fantasy Heist Ideas'
    'MRK Investors Have Opportunity to Lead Merck &amp; Co., Inc. Securities Fraud Lawsuit with the Schall Law Firm'
    'FTAI Aviation Ltd. (FTAI) Investors: March 18, 2025 Filing Deadline in Securities Class Action - Contact Kessler Topaz Meltzer & Check, LLP'
    'SIU investigates eThekwini Municipality’s Water and Sanitation Unit over alleged maladministration'
    '‘King of Crocodiles’ behind bars: Man faces criminal charges for brutal torture of crocodile'
    'Deliveroo driver arrested after returning to airport to collect bag containing €24,000 worth of pills'
    'Megha Engineering And Infrastructure Loses Rs 5.5 Crore In Phishing Scam'
    'Gachagua Speaks on Forgiving Ruto Amid Betrayal Allegations'
    'Pair stopped on suspicion of theft in town centre'
    'Bank manager grilled in Gono fraud case -Newsday Zimbabwe'
    'ROSEN, A LEADING LAW FIRM, Encourages Edison International Investors to Secure Counsel Before Important Deadline in Securities Class Action First Filed by the Firm – EIX'
    'Puspakom insists it’s compliant with JPJ inspection regs after accused of failing lorries unreasonably'
    'Food Head at FDA Resigns, Citing Numerous Job Cuts in Division'
    '18/02/2025: Sport: United malaise comes from the top and Amorim faces a battle to fix it'
    '18/02/2025: Sport: FA investigating alleged racist abuse of Mejbri'
    'West Tamar man charged with fraud after police seize four vehicles'
    'Cost claims by builders of new National Children’s Hospital rise another €100m amid fears final bill could still increase'
    'Drunk driver charged in resisting arrest'
    'The Great Social Security Heist: How Politicians Robbed You Blind and Blamed Everyone but…'
    'Former Yorkshire police officer jailed after sending inappropriate messages to crime victim'
    'Aboriginal heritage laws create chaos as first WA farmer fined for breach.'
    'Police arrest 2 suspected robbers in Benue'
    'Man accused of violent attack needs mental health assessed'
    'Rivervale man, 52, arrested after allegedly pointing laser at commercial planes and police helicopter'
    'Four NYC Deputy Mayors Resign Over Federal Involvement in Adams Case'
    'Akali Dal crisis deepens, Sikh body chief Harjinder Dhami resigns on moral grounds'
    'San Pedro Employee Charged with Theft of over $500,000'
    'Carla Freedman Concludes Her Service as United States Attorney'
    'Linklaters hires CMA litigation chief in London'
    'Man blinds pilots landing at Perth Airport with laser, then turns it on police chopper: police allege'
    'Woman’s Mugshot Breaks The Internet After She’s Arrested For Filling Her Ex-Boyfriend’s Gas Tank With Already Sucked Skittles & Coke Zero'
    'Alleged Leader of ''Zizians'' Group Arrested, Linked to Multiple US Killings'
    'Mizuho Employee Suspected of Theft in Latest Safety Box Case'
    'BORO PARK: NYPD and Shomrim Nab Mail Theft Suspects in Early Morning Arrest'
    'Questions to Ask a Debt Settlement Company'
    'Deputies chase, arrest ''highly intoxicated'' man on I-70 during ''holiday weekend chaos'''
    'Litigation in Corporate Law: Protecting Your Business from Legal Disputes'
    'How to Prevent Theft in the Workplace'
    'NT judge says ''tidal wave'' of arrests behind ''concerning'' custody waits'
    'I just played Killzone 2 at 60FPS and I feel like I should be arrested.'
    'Train crashes into truck in southern Oregon; 2 arrested'
    'Pinellas Park tile shop owner arrested again in scheme to defraud case'
    'Kessler Topaz Meltzer & Check, LLP - Class Action Announcement for Neumora Therapeutics, Inc. Investors: A Securities Fraud Class Action Lawsuit Was Filed Against Neumora Therapeutics, Inc.'
    'Petty theft at Hobby Lobby ends with arrest of three women'
    'Northglenn police searching for man accused of shoplifting, injuring officer'
    'Missouri Truck Company Owner Gets 9 Years For PPP Fraud, Other Felonies'
    'Man fined, vehicle seized for repeatedly colliding with another car'
    'Disney Forever Changes Toy Story Land, Corrects Wrongdoing'
    'Person in Crisis Dies Following Arrest by Quebec City Police (February 14-17, 2025)'
    'At least 3 NYC deputy mayors express intent to resign from Adams administration: Sources (Melissa Russo/NBC New York)'
    'Ravens kickers, Justin Tucker faces allegations of sexual misconduct from 16 massage therapist.'
    'Panini’s employee to be charged for involvement in theft and cashing of $7,000 check: Highland Heights Police Blotter'
    'New York city council calls for under-fire mayor''s resignation'
    'Pervasive Police Corruption in Albuquerque Explains Why a Teetotaler Was Arrested for DWI'
    'Honda allegedly to resume merger talks if Nissan''s Uchida leaves'
    'Boy, 10, in hospital with critical injuries after North Liverpool Rd crash'
    'the democrats should demand all contracts with a musk company be cancelled'
    'Eviction protection for those affected by wildfires up for LA County Board of Supervisors vote'
    'VA high school athletic director arrested while visiting family in Pennsylvania, accused of sending sexually explicit photos to underage student'
    'New York City Council Calls for Mayor Adams'' Resignation'
    'Road Workers Aid In Drunk Driver Arrest On SH 29'
    'Three held for theft, stolen property recovered'
    'Alaska Energy Metals Announces Resignation Of Director'
    'ROSEN, A LEADING INVESTOR RIGHTS LAW FIRM, Encourages Alarum Technologies Ltd. Investors to Secure Counsel Before Important Deadline in Securities Class Action – ALAR'
    'Man arrested after woman dies in hospital with ‘serious injuries’'
    'Suspect arrested in theft at construction site in The Villages'
    'Fact-Check: Did Washoe County mislead the legislature on homeless numbers? Yes.'
    'Does car insurance cover theft?'
)

$startRow = 195
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $titles[$i]
}
